# Update NATMI TPM-derived metrics for the Plg-F2r sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MuSCs -> ECs): receptor-expressing cell count increased from 2 to 3,
# which changes receptor average/total expression and all derived specificity
# and edge-weight columns that depend on it.
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.536281
$ws.Range("N2").Value = 13.608843
$ws.Range("O2").Value = 0.07292094190424039
$ws.Range("P2").Value = 0.07292094190424041
$ws.Range("Q2").Value = 0.08058249568399999
$ws.Range("R2").Value = 0.725242461156
$ws.Range("S2").Value = 0.07292094190424039
$ws.Range("T2").Value = 0.07292094190424041

# Row 3 (MuSCs -> FAPs): specificity columns renormalized because row 2's
# receptor expression total changed.
$ws.Range("O3").Value = 0.7154667412877611
$ws.Range("P3").Value = 0.7154667412877612
$ws.Range("S3").Value = 0.7154667412877611
$ws.Range("T3").Value = 0.7154667412877612

# Row 4 (MuSCs -> MuSCs): specificity columns renormalized as well.
$ws.Range("O4").Value = 0.2116123168079984
$ws.Range("P4").Value = 0.2116123168079984
$ws.Range("S4").Value = 0.2116123168079984
$ws.Range("T4").Value = 0.2116123168079984
